$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.994791030883789
$ws.Range("B1").Value = 2.362165927886963
$ws.Range("C1").Value = 2.478404760360718
$ws.Range("D1").Value = 3.148077011108398
$ws.Range("E1").Value = 2.23503303527832
